$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing existing rows 11-28 down to 12-29.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with this week's data (same market/category context,
# new observation date, price and related figures).
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44659
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112022
$ws.Range("G11").Value = "Arveja Verde"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 140
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24571
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Carahue"
$ws.Range("P11").Value = 983
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
